$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 23, shifting existing rows 23-26 down to 24-27.
$ws.Rows("23:23").Insert()

# Fill in the new row 23 with the new weekly data entry.
$ws.Range("A23").Value = 6
$ws.Range("B23").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C23").Value = "Metropolitana"
$ws.Range("D23").Value = 44644
$ws.Range("D23").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100101
$ws.Range("H23").Value = "Berries"
$ws.Range("I23").Value = 100101006
$ws.Range("J23").Value = "Higo"
$ws.Range("K23").Value = "Sin especificar"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 85
$ws.Range("N23").Value = 14000
$ws.Range("O23").Value = 14000
$ws.Range("P23").Value = 14000
$ws.Range("Q23").Value = "$/bandeja 7 kilos"
$ws.Range("R23").Value = "Región Metropolitana"
$ws.Range("S23").Value = 2000
$ws.Range("T23").Value = 7
